$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.749.95"
$ws.Range("E2").Value = "  +3.15%  "

$ws.Range("D3").Value = "2.519.06"
$ws.Range("E3").Value = "  +1.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").Value = "2.517.04"
$ws.Range("E9").Value = "  +1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.68%  "

$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "2.956.98"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.97"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "69.607.73"
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000178"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.89%  "

$ws.Range("D18").Value = "2.514.52"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.40"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.05"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.20%  "

$ws.Range("D28").Value = "2.643.22"
$ws.Range("E28").Value = "  +1.07%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "510.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.13%  "

$ws.Range("D31").Value = "0.0₃0889"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.96%  "

$ws.Range("E34").Value = "  +0.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.320"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.515"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0737"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").Value = "0.0₆0247"
$ws.Range("E51").Value = "  -2.75%  "

